$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.558.73'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.874.32'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("D4").Value = '''1.008'
$ws.Range("E4").Value = '  -1.63%  '
$ws.Range("D5").Value = '''315.83'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("D7").Value = '''0.5097'
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("E8").Value = '  -0.98%  '
$ws.Range("D9").Value = '''0.08358'
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = '''1.107'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").Value = '''41.86'
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").Value = '''6.227'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '1.883.42'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").Value = '''20.45'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '''7.287'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").Value = '''1.008'
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '''91.27'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Value = '''0.06738'
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").Value = '''17.74'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("D22").Value = '''5.915'
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").Value = '28.582.07'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '''11.15'
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("D25").Value = '''2.228'
$ws.Range("E25").Value = '  -1.60%  '
$ws.Range("D26").Value = '2.091.26'
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").Value = '''161.64'
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = '''20.65'
$ws.Range("D29").Value = '''2.420'
$ws.Range("E29").Value = '  +2.31%  '
$ws.Range("D30").Value = '''126.55'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").Value = '''0.1043'
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").Value = '''5.746'
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("D34").Value = '''3.608'
$ws.Range("E34").Value = '  -1.66%  '
$ws.Range("D35").Value = '''0.02452'
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("D36").Value = '''0.06564'
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("D37").Value = '''8.945'
$ws.Range("D38").Value = '''0.2167'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").Value = '''5.030'
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("D40").Value = '''1.183'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").Value = '''1.237'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").Value = '''0.6386'
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("D43").Value = '''11.12'
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("E44").Value = '  -1.52%  '
$ws.Range("D45").Value = '''0.6014'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").Value = '''13.08'
$ws.Range("E46").Value = '  +1.29%  '
$ws.Range("D47").Value = '''3.694'
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("D48").Value = '''2.007'
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = '''1.221'
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("D50").Value = '''122.11'
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("E51").Value = '  -10.39%  '
